# Sync attendance_reports from main repo.
#
# The "Recorded By" audit column (G) on the "Session Analysis Results"
# sheet lists the actor(s) that recorded/updated each attendance session,
# e.g. "dnasr281@gmail.com, System". Upstream normalized the ordering of
# these comma-separated actor lists so the canonical "System" actor is
# listed first (and, for the few rows that also carry a stray lowercase
# "system" duplicate, that duplicate now trails at the end instead of the
# front). Re-apply that same normalization here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact "before" -> "after" text replacements for the "Recorded By" column,
# derived from the upstream sync. Values not present in this map (e.g.
# "admin@admin.com, System") are left untouched.
$recordedByMap = @{
    "system, backup@backdoor.com, System" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Text
    if ($recordedByMap.ContainsKey($current)) {
        $cell.Value = $recordedByMap[$current]
    }
}
